$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "sex" column (D) previously had no data; every data row (2-6) now
# gets the value "U" in column D.
$ws.Range("D2:D6").Value = "U"

# Columns F ("sire") and G ("dam") are no longer populated for these rows -
# clear their contents so the cells disappear from the sheet entirely.
$ws.Range("F2:G6").ClearContents()

# The active selection moved from H6 to H9.
$ws.Range("H9").Select()
